$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header labels: "_old" -> "_FV2310", "_new" -> "_FV2404" ---
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value()
    if ($v -ne $null) {
        if ($v.EndsWith("_old")) {
            $cell.Value = $v.Substring(0, $v.Length - 4) + "_FV2310"
        } elseif ($v.EndsWith("_new")) {
            $cell.Value = $v.Substring(0, $v.Length - 4) + "_FV2404"
        }
    }
}

# --- 2) Turn the data range into an Excel Table (adds xl/tables/table1.xml
#        plus the <tableParts> wiring in the worksheet) ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U62"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = $null

# --- 3) Freeze the header row (pane split after row 1) ---
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$null
